# Natmi following Dr Hou advice:
# re-ran the LR-pair analysis with an extra "M1" target cluster, which
# updates the per-cluster stats on the existing rows and adds a new
# row for the "sCs" target cluster.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sending cluster FAPs -> Target cluster ECs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.5587383333333333
$ws.Range("H2").Value = 1.676215
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.042935333333334
$ws.Range("N2").Value = 9.128806
$ws.Range("O2").Value = 0.5356624298519768
$ws.Range("P2").Value = 0.5831288282815508
$ws.Range("Q2").Value = 1.700204616587778
$ws.Range("R2").Value = 15.30184154929
$ws.Range("S2").Value = 0.5356624298519768
$ws.Range("T2").Value = 0.5831288282815508

# Row 3 (Sending cluster FAPs -> Target cluster FAPs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.5587383333333333
$ws.Range("H3").Value = 1.676215
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.067915666666667
$ws.Range("N3").Value = 3.203747
$ws.Range("O3").Value = 0.1879902916822836
$ws.Range("P3").Value = 0.2046485853922773
$ws.Range("Q3").Value = 0.5966854197338889
$ws.Range("R3").Value = 5.370168777605
$ws.Range("S3").Value = 0.1879902916822836
$ws.Range("T3").Value = 0.2046485853922773

# Row 4 (Sending cluster FAPs -> Target cluster M1, relabeled from M2)
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.5587383333333333
$ws.Range("H4").Value = 1.676215
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02258333333333333
$ws.Range("N4").Value = 0.06775
$ws.Range("O4").Value = 0.003975451950942043
$ws.Range("P4").Value = 0.004327726771285868
$ws.Range("Q4").Value = 0.01261817402777778
$ws.Range("R4").Value = 0.11356356625
$ws.Range("S4").Value = 0.003975451950942043
$ws.Range("T4").Value = 0.004327726771285868

# Row 5 (Sending cluster FAPs -> Target cluster M2, relabeled from sCs)
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.5587383333333333
$ws.Range("H5").Value = 1.676215
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.160044
$ws.Range("N5").Value = 0.480132
$ws.Range("O5").Value = 0.02817330916767092
$ws.Range("P5").Value = 0.03066981712400039
$ws.Range("Q5").Value = 0.08942271782
$ws.Range("R5").Value = 0.80480446038
$ws.Range("S5").Value = 0.02817330916767092
$ws.Range("T5").Value = 0.03066981712400039

# New row 6 (new "sCs" target-cluster entry)
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Lama1"
$ws.Range("C6").Value = "Itga2"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5587383333333333
$ws.Range("H6").Value = 1.676215
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.3872175
$ws.Range("N6").Value = 2.774435
$ws.Range("O6").Value = 0.2441985173471266
$ws.Range("P6").Value = 0.1772250424308857
$ws.Range("Q6").Value = 0.7750915939208334
$ws.Range("R6").Value = 4.650549563525
$ws.Range("S6").Value = 0.2441985173471266
$ws.Range("T6").Value = 0.1772250424308857
